# This script applies the per-row crypto price/volume refresh captured in the
# authoritative OOXML diff (commit: "Updated cryptos list ... with GitHub Actions").
# It updates Price (D) / Volume(1h) (E) values, and for the two swapped row pairs
# (WrappedEther/WrappedBTC at rows 18-19, RenderToken/FirstDigitalUSD at rows 29-30,
# Maker/VeChain at rows 42-43) it also rewrites the Coin (B) and Link (C) columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new value, derived from the authoritative diff.
$updates = [ordered]@{
    'D2' = '63.084.20'
    'E2' = '  +6.33%  '
    'D3' = '3.117.40'
    'E3' = '  +3.86%  '
    'E4' = '  -0.02%  '
    'D5' = '586.49'
    'E5' = '  +3.93%  '
    'D6' = '143.56'
    'E6' = '  +3.33%  '
    'E7' = '  -0.04%  '
    'D8' = '3.105.90'
    'E8' = '  +3.91%  '
    'E9' = '  +2.53%  '
    'E10' = '  +9.17%  '
    'D11' = '5.75'
    'E11' = '  +10.49%  '
    'E12' = '  +2.35%  '
    'E13' = '  +5.36%  '
    'D14' = '35.55'
    'E14' = '  +5.25%  '
    'E15' = '  +0.80%  '
    'D16' = '3.631.64'
    'E16' = '  +3.72%  '
    'E17' = '  -1.14%  '
    'B18' = 'WrappedBTC'
    'C18' = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
    'D18' = '63.028.74'
    'E18' = '  +6.23%  '
    'B19' = 'WrappedEther'
    'C19' = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
    'D19' = '3.111.36'
    'E19' = '  +3.65%  '
    'D20' = '453.29'
    'E20' = '  +5.08%  '
    'E21' = '  +3.15%  '
    'E22' = '  +1.60%  '
    'D23' = '7.55'
    'E23' = '  +5.66%  '
    'E24' = '  +0.94%  '
    'D25' = '81.94'
    'E25' = '  +1.79%  '
    'E26' = '  +0.04%  '
    'E27' = '  +0.87%  '
    'E28' = '  +5.67%  '
    'B29' = 'FirstDigitalUSD'
    'C29' = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
    'D29' = '1.00'
    'E29' = '  +0.05%  '
    'B30' = 'RenderToken'
    'C30' = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
    'D30' = '8.28'
    'E30' = '  +4.95%  '
    'E31' = '  +11.92%  '
    'E32' = '  +12.48%  '
    'D33' = '27.15'
    'E33' = '  +5.32%  '
    'D34' = '1.04'
    'E34' = '  +4.18%  '
    'D35' = '0.0₃0806'
    'E35' = '  +6.35%  '
    'D36' = '2.32'
    'E36' = '  +9.69%  '
    'E37' = '  +1.33%  '
    'E38' = '  +3.89%  '
    'D39' = '3.03'
    'E39' = '  +10.11%  '
    'D40' = '8.78'
    'D41' = '427.06'
    'E41' = '  +4.18%  '
    'B42' = 'VeChain'
    'C42' = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
    'D42' = '0.0375'
    'E42' = '  +5.60%  '
    'B43' = 'Maker'
    'C43' = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
    'D43' = '2.949.32'
    'E43' = '  +6.23%  '
    'E44' = '  +9.74%  '
    'D45' = '0.112'
    'E45' = '  +3.60%  '
    'E46' = '  +7.57%  '
    'D47' = '125.89'
    'E47' = '  +1.87%  '
    'D49' = '34.80'
    'E49' = '  -0.51%  '
    'E50' = '  +1.16%  '
    'D51' = '24.79'
    'E51' = '  +5.19%  '
}

# Cells whose new value is a plain decimal number (e.g. "1.00", "34.80") need to be
# explicitly formatted as text first, otherwise Excel would silently coerce them to
# a number and drop significant trailing zeros (e.g. "1.00" -> 1, "34.80" -> 34.8).
$forceText = @('D5', 'D6', 'D11', 'D14', 'D20', 'D23', 'D25', 'D29', 'D30', 'D33', 'D34', 'D36', 'D39', 'D40', 'D41', 'D42', 'D45', 'D47', 'D49', 'D51')

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    if ($forceText -contains $ref) {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $updates[$ref]
}
